# Apply the "update initial pieces and placements" edit:
#  - Rename the "CAP" piece label to "CPT" (cells C3 and C9)
#  - Update the active selection to I15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the piece label text from "CAP" to "CPT"
$ws.Range("C3").Value = "CPT"
$ws.Range("C9").Value = "CPT"

# Update the selected cell/range shown when the sheet is active
$ws.Activate()
$ws.Range("I15").Select()
